# Auto-applies the crypto price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.604.44"
$ws.Range("E2").Value = "  +7.32%  "
$ws.Range("D3").Value = "3.625.33"
$ws.Range("E3").Value = "  +7.01%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'593.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.17%  "
$ws.Range("D6").Value = "'191.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.56%  "
$ws.Range("D7").Value = "'0.649"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.99%  "
$ws.Range("D8").Value = "3.606.60"
$ws.Range("E8").Value = "  +6.71%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "'0.180"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.03%  "
$ws.Range("D11").Value = "'0.663"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.56%  "
$ws.Range("D12").Value = "'57.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.60%  "
$ws.Range("D13").Value = "'0.0000295"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.56%  "
$ws.Range("D14").Value = "'9.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.61%  "
$ws.Range("D15").Value = "4.208.42"
$ws.Range("E15").Value = "  +7.08%  "
$ws.Range("D16").Value = "3.632.58"
$ws.Range("E16").Value = "  +7.69%  "
$ws.Range("D17").Value = "'19.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.47%  "
$ws.Range("D18").Value = "70.619.43"
$ws.Range("E18").Value = "  +7.31%  "
$ws.Range("D19").Value = "'12.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.12%  "
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("E21").Value = "  +5.43%  "
$ws.Range("D22").Value = "'494.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.89%  "
$ws.Range("D23").Value = "'5.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.15%  "
$ws.Range("D24").Value = "'16.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +15.85%  "
$ws.Range("E25").Value = "  +9.07%  "
$ws.Range("D26").Value = "'91.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.72%  "
$ws.Range("D27").Value = "'3.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.24%  "
$ws.Range("D28").Value = "'11.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.58%  "
$ws.Range("D29").Value = "'9.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.04%  "
$ws.Range("D30").Value = "'32.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.72%  "
$ws.Range("D31").Value = "'7.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.47%  "
$ws.Range("D32").Value = "'12.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.77%  "
$ws.Range("D33").Value = "'618.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.49%  "
$ws.Range("E34").Value = "  +8.49%  "
$ws.Range("D35").Value = "'65.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.91%  "
$ws.Range("D36").Value = "0.0₃0832"
$ws.Range("E36").Value = "  +11.31%  "
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").Value = "'0.406"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.70%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.148"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.87%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "'38.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.60%  "
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("D42").Value = "3.346.86"
$ws.Range("E42").Value = "  +7.60%  "
$ws.Range("D43").Value = "'3.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.46%  "
$ws.Range("D44").Value = "'0.0449"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.54%  "
$ws.Range("D45").Value = "'2.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.42%  "
$ws.Range("D46").Value = "'3.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.18%  "
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("D48").Value = "'9.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.65%  "
$ws.Range("D49").Value = "'2.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.35%  "
$ws.Range("D50").Value = "'3.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.16%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.01%  "
